# Apply the changes described in the commit to the workbook.
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("India Data")

# E15: drop the formula, keep the literal value 0.15
$wsData.Range("E15").Value = 0.15

# D16: formula changes to reference D9 instead of the SYVbT-freight ratio
$wsData.Range("D16").Formula = "=D9"

# E19: drop the formula, becomes a literal 0
$wsData.Range("E19").Value = 0

# E26: drop the formula, becomes a literal 0
$wsData.Range("E26").Value = 0

# E29: drop the formula, keep the literal value 0.15
$wsData.Range("E29").Value = 0.15

# E33: drop the formula, becomes a literal 0
$wsData.Range("E33").Value = 0

# E85: drop the formula, set literal value to 0.3 (was 0.15)
$wsData.Range("E85").Value = 0.3

# E89 keeps its formula (=E85*'India Assumptions'!$A$36); it will recalc automatically

# Sheet view changes: India Data becomes the active/selected sheet
$wsAbout = $wb.Worksheets.Item("About")

$wsData.Activate()
$wsData.Application.ActiveWindow.ScrollRow = 21
$wsData.Range("E34").Select()
